# Issue #13: allow metadata files to have a dedicated row of machine
# ("slug"-style) column identifiers, so two columns can later be related
# to build hierarchical SKOS concepts.
#
# Concretely: insert a brand-new row directly under the human-readable
# header row (row 1). This new row 2 carries the machine name for each
# column. Every row that used to be row 2/3/4 (the semantic-vocabulary /
# dim-type / datatype rows) shifts down one position to rows 3/4/5. The
# old, orphaned row 5 (which only held a stray "mapping-ano.xlsx" value in
# column H) is dropped once the real datatype row lands on row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2..5 down to 3..6 and open up a blank row 2.
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with the machine-readable column identifiers, in the
# same left-to-right column order as the header row.
$ws.Range("A2").Value = "poblacion"
$ws.Range("B2").Value = "ccaa-nombre"
$ws.Range("C2").Value = "comarca-nombre"
$ws.Range("D2").Value = "comarca-codigo"
$ws.Range("E2").Value = "provincia-codigo"
$ws.Range("F2").Value = "municipio-codigo"
$ws.Range("G2").Value = "provincia-nombre"
$ws.Range("H2").Value = "ano"
$ws.Range("I2").Value = "municipio-nombre"

# The old row 5 (now pushed to row 6) only ever had a single leftover cell
# (H6 = "mapping-ano.xlsx"); the real datatype row has already taken its
# rightful place at row 5, so this now-empty-ish leftover row is removed.
$ws.Rows.Item(6).Delete()
